$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is numeric-looking text (e.g. "1.006")
# must be forced to Text format first, otherwise Excel auto-converts the
# assigned string into a Number, same as the source data (stored as text).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.809.64"
$ws.Range("E2").Value = "  +1.60%  "

$ws.Range("D3").Value = "1.887.64"
$ws.Range("E3").Value = "  +1.60%  "

$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.47%  "

$ws.Range("D5").Value = "333.36"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("D7").Value = "0.4719"
$ws.Range("E7").Value = "  +3.16%  "

$ws.Range("D8").Value = "0.3945"
$ws.Range("E8").Value = "  +0.44%  "

$ws.Range("D9").Value = "47.77"
$ws.Range("E9").Value = "  +1.40%  "

$ws.Range("D10").Value = "0.08079"
$ws.Range("E10").Value = "  +1.74%  "

$ws.Range("D11").Value = "1.028"
$ws.Range("E11").Value = "  +1.52%  "

$ws.Range("D12").Value = "22.26"
$ws.Range("E12").Value = "  +3.70%  "

$ws.Range("D13").Value = "1.887.62"
$ws.Range("E13").Value = "  +1.97%  "

$ws.Range("D14").Value = "5.989"
$ws.Range("E14").Value = "  +1.17%  "

$ws.Range("D15").Value = "7.146"
$ws.Range("E15").Value = "  -0.07%  "

$ws.Range("D16").Value = "1.009"
$ws.Range("E16").Value = "  +0.60%  "

$ws.Range("D17").Value = "0.06760"
$ws.Range("E17").Value = "  +2.23%  "

$ws.Range("D18").Value = "87.40"
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("E19").Value = "  +1.74%  "

$ws.Range("D20").Value = "17.37"
$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("D21").Value = "1.006"
$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "27.823.45"
$ws.Range("E22").Value = "  +1.66%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.545"
$ws.Range("E23").Value = "  +1.05%  "

$ws.Range("D24").Value = "11.04"
$ws.Range("E24").Value = "  +1.03%  "

$ws.Range("D25").Value = "2.319"
$ws.Range("E25").Value = "  +0.78%  "

$ws.Range("D26").Value = "2.101.86"
$ws.Range("E26").Value = "  +1.38%  "

$ws.Range("D27").Value = "159.32"
$ws.Range("E27").Value = "  +3.69%  "

$ws.Range("D28").Value = "20.17"
$ws.Range("E28").Value = "  +0.46%  "

$ws.Range("D29").Value = "2.117"
$ws.Range("E29").Value = "  +2.26%  "

$ws.Range("D30").Value = "5.600"
$ws.Range("E30").Value = "  +2.41%  "

$ws.Range("D31").Value = "122.14"
$ws.Range("E31").Value = "  +0.37%  "

$ws.Range("D32").Value = "0.9878"
$ws.Range("E32").Value = "  +3.68%  "

$ws.Range("D33").Value = "0.09493"
$ws.Range("E33").Value = "  +1.05%  "

$ws.Range("D34").Value = "1.457"
$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").Value = "3.614"
$ws.Range("E35").Value = "  +0.64%  "

$ws.Range("D36").Value = "5.370"
$ws.Range("E36").Value = "  +1.87%  "

$ws.Range("D37").Value = "0.06160"
$ws.Range("E37").Value = "  +2.15%  "

$ws.Range("D38").Value = "0.02274"
$ws.Range("E38").Value = "  +1.99%  "

$ws.Range("E39").Value = "  +1.03%  "

$ws.Range("D40").Value = "8.116"
$ws.Range("E40").Value = "  +0.79%  "

$ws.Range("D41").Value = "0.6017"
$ws.Range("E41").Value = "  +1.60%  "

$ws.Range("D42").Value = "0.1897"
$ws.Range("E42").Value = "  +0.60%  "

$ws.Range("D43").Value = "10.34"
$ws.Range("E43").Value = "  +1.79%  "

$ws.Range("D44").Value = "1.261"
$ws.Range("E44").Value = "  -1.61%  "

$ws.Range("D45").Value = "0.5736"
$ws.Range("E45").Value = "  +2.08%  "

$ws.Range("D46").Value = "12.22"
$ws.Range("E46").Value = "  +0.86%  "

$ws.Range("D47").Value = "1.955"
$ws.Range("E47").Value = "  +1.96%  "

$ws.Range("D48").Value = "3.394"
$ws.Range("E48").Value = "  -0.06%  "

$ws.Range("E49").Value = "  +2.42%  "

$ws.Range("D50").Value = "113.67"
$ws.Range("E50").Value = "  +4.99%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.00000000302"
$ws.Range("E51").Value = "  +7.74%  "
